$d = $word.ActiveDocument

# Anchor on the "LOB1019: Física II (Requisito fraco)" paragraph, which stays
# in place. The three paragraphs that directly follow it -- a blank paragraph,
# "Ver no Jupiter Salvar em pdf Salvar em docx", and the "© 2020 ..." footer
# line -- must be removed, while the paragraphs that come after them (the
# trailing blank paragraph and the page-break paragraph) are left untouched.
$rng = $d.Content
$found = $rng.Find.Execute(
    "LOB1019: Física II (Requisito fraco)",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $anchorPara = $rng.Paragraphs(1)

    $blankPara = $anchorPara.Next()
    $jupiterPara = $blankPara.Next()
    $copyrightPara = $jupiterPara.Next()

    $deleteRange = $d.Range($blankPara.Range.Start, $copyrightPara.Range.End)
    $deleteRange.Delete()
}
